# Update countries & provincias Spain
# Refreshes the COVID-19 country table on sheet "Pais":
#   - updates the "last updated" timestamp in A1
#   - refreshes Casos totales/Nuevos casos/Casos activos/Recuperados/
#     Casos criticos/Muertes hoy/Muertes for the countries whose figures
#     changed between the 13:52 and 14:22 snapshots
#   - re-sorts a handful of same-total-cases countries that swapped
#     positions in the (descending, by Casos totales) ranking

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Datos actualizados a ..." banner above the table
$ws.Range("A1").Value = "Datos actualizados a 5 de Abril de 2020 a las 14:22"

# Rows (A:H = Pais, Casos totales, Nuevos casos, Casos activos, Recuperados,
# Casos criticos, Muertes hoy, Muertes) that changed, keyed by row number.
$rowUpdates = @(
    @{ Row = 13; Values = @("Suiza", 21100, 595, 6415, 14005, 391, 14, 680) },
    @{ Row = 15; Values = @("Paises Bajos", 17851, 1224, 250, 15835, 1360, 115, 1766) },
    @{ Row = 17; Values = @("Austria", 11897, 116, 2998, 8695, 244, 18, 204) },
    @{ Row = 22; Values = @("Suecia", 6830, 387, 205, 6224, 541, 28, 401) },
    @{ Row = 58; Values = @("Croacia", 1182, 56, 119, 1047, 39, 4, 16) },
    @{ Row = 63; Values = @("Marruecos", 961, 42, 70, 822, 1, 10, 69) },
    @{ Row = 78; Values = @("Republica de Macedonia", 555, 72, 23, 514, 15, 1, 18) },
    @{ Row = 79; Values = @("Tunez", 553, 0, 5, 529, 39, 1, 19) },
    @{ Row = 80; Values = @("Letonia", 533, 24, 1, 531, 4, 0, 1) },
    @{ Row = 81; Values = @("Libano", 527, 7, 54, 455, 26, 1, 18) },
    @{ Row = 82; Values = @("Bulgaria", 522, 19, 37, 467, 26, 1, 18) },
    @{ Row = 83; Values = @("Principado de Andorra", 501, 35, 26, 457, 12, 1, 18) },
    @{ Row = 84; Values = @("Eslovaquia", 485, 14, 10, 474, 3, 0, 1) },
    @{ Row = 109; Values = @("Sri Lanka", 171, 5, 29, 137, 5, 0, 5) },
    @{ Row = 146; Values = @("Puerto Rico", 39, 0, 1, 36, 0, 0, 2) },
    @{ Row = 147; Values = @("Zambia", 39, 0, 3, 35, 0, 0, 1) },
    @{ Row = 156; Values = @("Birmania", 21, 0, 0, 20, 0, 0, 1) },
    @{ Row = 157; Values = @("Haiti", 21, 1, 1, 20, 0, 0, 0) },
    @{ Row = 173; Values = @("Fiyi", 12, 0, 0, 12, 0, 0, 0) },
    @{ Row = 174; Values = @("Granada", 12, 0, 0, 12, 2, 0, 0) },
    @{ Row = 181; Values = @("Angola", 10, 0, 2, 6, 0, 0, 2) },
    @{ Row = 182; Values = @("Sudan", 10, 0, 2, 6, 0, 0, 2) },
    @{ Row = 184; Values = @("Republica del Chad", 9, 0, 0, 9, 0, 0, 0) },
    @{ Row = 185; Values = @("Suazilandia", 9, 0, 0, 9, 0, 0, 0) },
    @{ Row = 187; Values = @("Nepal", 9, 0, 1, 8, 0, 0, 0) },
    @{ Row = 188; Values = @("Zimbabue", 9, 0, 0, 8, 0, 0, 1) },
    @{ Row = 192; Values = @("San Vicente y las Granadinas", 7, 0, 1, 6, 0, 0, 0) },
    @{ Row = 193; Values = @("Cabo Verde", 7, 0, 0, 6, 0, 0, 1) },
    @{ Row = 201; Values = @("Belice", 4, 0, 0, 4, 0, 0, 0) },
    @{ Row = 202; Values = @("Malaui", 4, 0, 0, 4, 0, 0, 0) },
    @{ Row = 206; Values = @("Burundi", 3, 0, 0, 3, 0, 0, 0) },
    @{ Row = 208; Values = @("Islas Virgenes Britanicas", 3, 0, 0, 3, 0, 0, 0) },
    @{ Row = 210; Values = @("Sudan del Sur", 1, 1, 0, 1, 0, 0, 0) },
    @{ Row = 211; Values = @("Islas Malvinas", 1, 0, 0, 1, 0, 0, 0) },
    @{ Row = 212; Values = @("Timor Oriental", 1, 0, 0, 1, 0, 0, 0) },
    @{ Row = 213; Values = @("Papua Nueva Guinea", 1, 0, 0, 1, 0, 0, 0) },
    @{ Row = 214; Values = @("San Pedro y Miquelon", 1, 1, 0, 1, 0, 0, 0) }
)

foreach ($u in $rowUpdates) {
    $r = $u.Row
    $vals = $u.Values
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = $i + 1
        $ws.Cells.Item($r, $col).Value = $vals[$i]
    }
}
